$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# B2: "moloko" -> "NEWMOLOKO"
$ws.Range("B2").Value = "NEWMOLOKO"

# E3: "true" -> "false" (copy the existing text "false" from E4 so it stays
# a text/shared-string value instead of being auto-coerced to a Boolean)
$ws.Range("E4").Copy()
$ws.Range("E3").PasteSpecial(-4163)

# A5: 3 -> 4
$ws.Range("A5").Value = 4

# B5: "cucus" -> "cucusNEW"
$ws.Range("B5").Value = "cucusNEW"

# Update the active selection to J6
$ws.Range("J6").Select()
